$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Add the two new ISMIP6 initMIP mapped-collection folder rows (6 & 7)
# Reuse the formatting already used by the existing data rows so no new
# fonts/number-formats are introduced.
# ---------------------------------------------------------------------

# A-column style (folder path cells)
$ws.Range("A4").Copy() | Out-Null
$ws.Range("A6:A7").PasteSpecial(-4122) | Out-Null

# B-column style (modeling groups, Menlo font / reading order)
$ws.Range("B4").Copy() | Out-Null
$ws.Range("B6:B7").PasteSpecial(-4122) | Out-Null

# C-column style for the new rows matches the plain bordered style (same as column A)
$ws.Range("A4").Copy() | Out-Null
$ws.Range("C6:C7").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# New values: ISMIP6 initMIP mapped collection folders
$ws.Range("A6").Value = "/projects/grid/ghub/ISMIP6/initMIP/AIS/output_original"
$ws.Range("B6").Value = "ARC,AWI,CPOM,DMI,DOE,IGE,ILTS,IMAU,JPL1,LSCE,NCAR,PIK,PSU,UCIJPL,ULB,VUB"
$ws.Range("C6").Value = "Antarctica"

$ws.Range("A7").Value = "/projects/grid/ghub/ISMIP6/initMIP/GrIS/output_original"
$ws.Range("B7").Value = "ARC,AWI,BGC,DMI,ILTS,ILTS_PIK,IMAU,ISMIP6,JPL1,LANL,LGGE,LSCE,MIROC,MPIM,UAF,UCIJPL,ULB,VUB"
$ws.Range("C7").Value = "Greenland"

# ---------------------------------------------------------------------
# Add a decorative empty bordered "box" spanning A8:E10 (white fill,
# thin gray outside border only, no interior gridlines)
# ---------------------------------------------------------------------
$gray  = 11184810
$white = 16777215

function SetEdge($r, $item) {
    $r.Borders.Item($item).Color = $gray
    $r.Borders.Item($item).LineStyle = 1
    $r.Borders.Item($item).Weight = 2
}

# Row 8 - top edge of the box
$r = $ws.Range("A8");    $r.Interior.Color = $white; SetEdge $r 7; SetEdge $r 8
$r = $ws.Range("B8:D8"); $r.Interior.Color = $white; SetEdge $r 8
$r = $ws.Range("E8");    $r.Interior.Color = $white; SetEdge $r 8; SetEdge $r 10

# Row 9 - middle of the box (no top/bottom border)
$r = $ws.Range("A9");    $r.Interior.Color = $white; SetEdge $r 7
$r = $ws.Range("B9:D9"); $r.Interior.Color = $white
$r = $ws.Range("E9");    $r.Interior.Color = $white; SetEdge $r 10

# Row 10 - bottom edge of the box
$r = $ws.Range("A10");    $r.Interior.Color = $white; SetEdge $r 7; SetEdge $r 9
$r = $ws.Range("B10:D10"); $r.Interior.Color = $white; SetEdge $r 9
$r = $ws.Range("E10");    $r.Interior.Color = $white; SetEdge $r 9; SetEdge $r 10
